$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Bold "2.2" inside " que foi explicado na seção 2.2, o que acabou..."
# ---------------------------------------------------------------------------
$ctx1 = $d.Content.Duplicate
$found1 = $ctx1.Find.Execute("na seção 2.2, o que acabou", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $inner1 = $d.Range($ctx1.Start, $ctx1.End)
    $found1b = $inner1.Find.Execute("2.2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found1b) {
        $inner1.Font.Bold = 1
    }
}

# ---------------------------------------------------------------------------
# 2) Remove the "_GoBack" bookmark that currently splits
#    " algum avião de acordo com as maiores filas, para t" / "entar manter..."
#    into two runs, and merge that text back into a single run.
# ---------------------------------------------------------------------------
$pre2 = $d.Content.Duplicate
$foundPre2 = $pre2.Find.Execute(" algum avião de acordo com as maiores filas, para t", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$post2 = $d.Content.Duplicate
$foundPost2 = $post2.Find.Execute("entar manter um padrão de tamanho das filas. E com isso era obtido um resultado insatisfatório de aviões caídos, conforme tabela abaixo: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPre2 -and $foundPost2) {
    $whole2 = $d.Range($pre2.Start, $post2.End)
    # Force a genuine content change so the engine regenerates the run (and
    # drops the now-interior bookmark) instead of treating it as a no-op.
    $whole2.Text = "##TEMP_PLACEHOLDER##"
    $reset2 = $d.Range($pre2.Start, $pre2.Start + "##TEMP_PLACEHOLDER##".Length)
    $reset2.Text = " algum avião de acordo com as maiores filas, para tentar manter um padrão de tamanho das filas. E com isso era obtido um resultado insatisfatório de aviões caídos, conforme tabela abaixo: "
}

# ---------------------------------------------------------------------------
# 3) Bold "2.2" inside "para deixa-lo da forma que é explicado na seção 2.2"
# ---------------------------------------------------------------------------
$ctx3 = $d.Content.Duplicate
$found3 = $ctx3.Find.Execute("explicado na seção 2.2, deixou-se", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $inner3 = $d.Range($ctx3.Start, $ctx3.End)
    $found3b = $inner3.Find.Execute("2.2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found3b) {
        $inner3.Font.Bold = 1
    }
}

# ---------------------------------------------------------------------------
# 4) Add a "_GoBack" bookmark right after "...Conforme tabela abaixo:" in the
#    paragraph that ends "...passou a não cair nenhum avião mais. Conforme
#    tabela abaixo:" (the one following the ProcessarIteracao refactor text).
#    A zero-length Range confuses bookmark placement in this host, so we
#    anchor on a throwaway character, wrap the bookmark around it, then
#    delete the character - the (now zero-width) bookmark stays put.
# ---------------------------------------------------------------------------
$ctx4 = $d.Content.Duplicate
$found4 = $ctx4.Find.Execute("não cair nenhum avião mais. Conforme tabela abaixo:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $anchorPos = $ctx4.End
    $ctx4.InsertAfter("Q")
    $anchorRange = $d.Range($anchorPos, $anchorPos + 1)
    $d.Bookmarks.Add("_GoBack", $anchorRange)
    $cleanupRange = $d.Range($anchorPos, $anchorPos + 1)
    $cleanupRange.Text = ""
}
